# Update loading_percent results (columns C:K, rows 2:25) for the
# "case with 380 kV done" run. Column I stays 0 (unchanged).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object "object[,]" 24,9
$arr[0,0] = 4.914374398633892
$arr[0,1] = 7.947514790248945
$arr[0,2] = 13.174742132032
$arr[0,3] = 38.767471982248
$arr[0,4] = 44.79128294397915
$arr[0,5] = 18.12062112128907
$arr[0,6] = 0
$arr[0,7] = 10.14078818510916
$arr[0,8] = 19.01930553588639
$arr[1,0] = 4.752283267633945
$arr[1,1] = 7.896805767506206
$arr[1,2] = 13.12060197606458
$arr[1,3] = 38.81922246966382
$arr[1,4] = 44.82755715438158
$arr[1,5] = 18.19467524776741
$arr[1,6] = 0
$arr[1,7] = 10.152379113861
$arr[1,8] = 18.47782554612911
$arr[2,0] = 4.651546961308906
$arr[2,1] = 7.866196549291487
$arr[2,2] = 13.08983786114731
$arr[2,3] = 38.864911999623
$arr[2,4] = 44.87057412266776
$arr[2,5] = 18.24513601421002
$arr[2,6] = 0
$arr[2,7] = 10.16134422693544
$arr[2,8] = 18.14020720115404
$arr[3,0] = 4.610264573821107
$arr[3,1] = 7.853863581349735
$arr[3,2] = 13.07793262852692
$arr[3,3] = 38.8870086403184
$arr[3,4] = 44.89328051290686
$arr[3,5] = 18.26694769475153
$arr[3,6] = 0
$arr[3,7] = 10.16546172490762
$arr[3,8] = 18.001537295394
$arr[4,0] = 4.603397918857534
$arr[4,1] = 7.851824424944985
$arr[4,2] = 13.07599413442354
$arr[4,3] = 38.89088715909174
$arr[4,4] = 44.89736228544169
$arr[4,5] = 18.27064472062479
$arr[4,6] = 0
$arr[4,7] = 10.16617344526376
$arr[4,8] = 17.97845185869432
$arr[5,0] = 4.650991049560919
$arr[5,1] = 7.866029642853883
$arr[5,2] = 13.08967473618697
$arr[5,3] = 38.86519595102842
$arr[5,4] = 44.87085944563936
$arr[5,5] = 18.24542512786603
$arr[5,6] = 0
$arr[5,7] = 10.16139787856122
$arr[5,8] = 18.13834117188807
$arr[6,0] = 4.858779156810874
$arr[6,1] = 7.929926045837095
$arr[6,2] = 13.15556542900794
$arr[6,3] = 38.78241555665657
$arr[6,4] = 44.79946246730119
$arr[6,5] = 18.14511512285735
$arr[6,6] = 0
$arr[6,7] = 10.14440079792926
$arr[6,8] = 18.83379156042855
$arr[7,0] = 5.253691665083892
$arr[7,1] = 8.059036481984332
$arr[7,2] = 13.30404566854689
$arr[7,3] = 38.73132801624955
$arr[7,4] = 44.8255683613028
$arr[7,5] = 17.9883015070056
$arr[7,6] = 0
$arr[7,7] = 10.12575849490272
$arr[7,8] = 20.14832106888554
$arr[8,0] = 5.532625868429681
$arr[8,1] = 8.155717858429089
$arr[8,2] = 13.42430830883723
$arr[8,3] = 38.76259916576542
$arr[8,4] = 44.94769917870005
$arr[8,5] = 17.89780368532932
$arr[8,6] = 0
$arr[8,7] = 10.12104454207346
$arr[8,8] = 21.07382084891195
$arr[9,0] = 5.656448000749541
$arr[9,1] = 8.199988515773597
$arr[9,2] = 13.48130983141607
$arr[9,3] = 38.7919226883493
$arr[9,4] = 45.02584792182532
$arr[9,5] = 17.86208217821297
$arr[9,6] = 0
$arr[9,7] = 10.12085407813752
$arr[9,8] = 21.48422390287061
$arr[10,0] = 5.702849411337894
$arr[10,1] = 8.216785079523611
$arr[10,2] = 13.50321224594354
$arr[10,3] = 38.80520623434451
$arr[10,4] = 45.05869860830145
$arr[10,5] = 17.84934513312908
$arr[10,6] = 0
$arr[10,7] = 10.12106296647934
$arr[10,8] = 21.63796907814046
$arr[11,0] = 5.692878396274287
$arr[11,1] = 8.213166356294115
$arr[11,2] = 13.49848126906081
$arr[11,3] = 38.80224835592267
$arr[11,4] = 45.05147867575287
$arr[11,5] = 17.85205304642086
$arr[11,6] = 0
$arr[11,7] = 10.12100548131439
$arr[11,8] = 21.60493345799927
$arr[12,0] = 5.660275487286752
$arr[12,1] = 8.201369787482077
$arr[12,2] = 13.48310547369811
$arr[12,3] = 38.79297181876665
$arr[12,4] = 45.02848525524552
$arr[12,5] = 17.86101842477414
$arr[12,6] = 0
$arr[12,7] = 10.12086563184188
$arr[12,8] = 21.49690674239621
$arr[13,0] = 5.64024050863063
$arr[13,1] = 8.194147964641768
$arr[13,2] = 13.47372827723907
$arr[13,3] = 38.7875736861485
$arr[13,4] = 45.01482549228703
$arr[13,5] = 17.86661304524489
$arr[13,6] = 0
$arr[13,7] = 10.12081656505442
$arr[13,8] = 21.43051644500734
$arr[14,0] = 5.524468012057516
$arr[14,1] = 8.15282969146482
$arr[14,2] = 13.42062820854353
$arr[14,3] = 38.76098738976163
$arr[14,4] = 44.94304746963446
$arr[14,5] = 17.90024834380861
$arr[14,6] = 0
$arr[14,7] = 10.12109631394284
$arr[14,8] = 21.04677426346608
$arr[15,0] = 5.452625282199937
$arr[15,1] = 8.127549782873881
$arr[15,2] = 13.38863163328007
$arr[15,3] = 38.74855133924987
$arr[15,4] = 44.90480842360486
$arr[15,5] = 17.9222822968123
$arr[15,6] = 0
$arr[15,7] = 10.12176845111241
$arr[15,8] = 20.80853927543847
$arr[16,0] = 5.411017050743381
$arr[16,1] = 8.113037294449576
$arr[16,2] = 13.37044470991657
$arr[16,3] = 38.7428192436333
$arr[16,4] = 44.88494056629839
$arr[16,5] = 17.93546789172781
$arr[16,6] = 0
$arr[16,7] = 10.12233896655937
$arr[16,8] = 20.67052146768999
$arr[17,0] = 5.396881563383083
$arr[17,1] = 8.108128685389037
$arr[17,2] = 13.36432451759952
$arr[17,3] = 38.74112217791173
$arr[17,4] = 44.87857845554043
$arr[17,5] = 17.94002007927481
$arr[17,6] = 0
$arr[17,7] = 10.1225637183139
$arr[17,8] = 20.62362522623256
$arr[18,0] = 5.460303035792886
$arr[18,1] = 8.130238058657387
$arr[18,2] = 13.3920153895428
$arr[18,3] = 38.74972806297396
$arr[18,4] = 44.90865889401325
$arr[18,5] = 17.91988368023272
$arr[18,6] = 0
$arr[18,7] = 10.12167786527288
$arr[18,8] = 20.83400345213321
$arr[19,0] = 5.669865322590638
$arr[19,1] = 8.20483392644516
$arr[19,2] = 13.48761321260419
$arr[19,3] = 38.79563736352332
$arr[19,4] = 45.03515053480248
$arr[19,5] = 17.85836358487358
$arr[19,6] = 0
$arr[19,7] = 10.12089908276393
$arr[19,8] = 21.52868305510381
$arr[20,0] = 5.803967822381438
$arr[20,1] = 8.253770093410713
$arr[20,2] = 13.55193444922832
$arr[20,3] = 38.83834655655484
$arr[20,4] = 45.13680663551277
$arr[20,5] = 17.82276388787403
$arr[20,6] = 0
$arr[20,7] = 10.12202798496031
$arr[20,8] = 21.97293295909242
$arr[21,0] = 5.73267036687334
$arr[21,1] = 8.22763827759376
$arr[21,2] = 13.5174407025004
$arr[21,3] = 38.81438737374423
$arr[21,4] = 45.08081231915688
$arr[21,5] = 17.84134036852799
$arr[21,6] = 0
$arr[21,7] = 10.12127562403788
$arr[21,8] = 21.73676440092207
$arr[22,0] = 5.456832869520132
$arr[22,1] = 8.129022622230485
$arr[22,2] = 13.39048494347775
$arr[22,3] = 38.749191651236
$arr[22,4] = 44.90691150505209
$arr[22,5] = 17.92096648178196
$arr[22,6] = 0
$arr[22,7] = 10.12171824570635
$arr[22,8] = 20.82249436801428
$arr[23,0] = 5.148590472141603
$arr[23,1] = 8.023751803946398
$arr[23,2] = 13.26187139830134
$arr[23,3] = 38.73312699033455
$arr[23,4] = 44.80052525000801
$arr[23,5] = 18.02641186627226
$arr[23,6] = 0
$arr[23,7] = 10.12922582631263
$arr[23,8] = 19.79909353781898
$ws.Range("C2:K25").Value = $arr
